$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.813.44"
$ws.Range("E2").Value = "  -3.03%  "
$ws.Range("D3").Value = "1.615.46"
$ws.Range("E3").Value = "  -3.44%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.71"
$ws.Range("E5").Value = "  -2.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3928"
$ws.Range("E7").Value = "  -0.70%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3832"
$ws.Range("E8").Value = "  -2.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.002"
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "49.31"
$ws.Range("E10").Value = "  -1.95%  "
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.352"
$ws.Range("E11").Value = "  -3.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08428"
$ws.Range("E12").Value = "  -2.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.60"
$ws.Range("E13").Value = "  -6.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.032"
$ws.Range("E14").Value = "  -3.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.558"
$ws.Range("E15").Value = "  -1.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001276"
$ws.Range("E16").Value = "  -2.88%  "
$ws.Range("D17").Value = "1.619.31"
$ws.Range("E17").Value = "  -3.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.71"
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06927"
$ws.Range("E19").Value = "  -1.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.95"
$ws.Range("E20").Value = "  -5.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.798"
$ws.Range("E21").Value = "  -3.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.41"
$ws.Range("E23").Value = "  -3.38%  "
$ws.Range("D24").Value = "23.818.64"
$ws.Range("E24").Value = "  -3.00%  "
$ws.Range("E25").Value = "  +4.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.838"
$ws.Range("E26").Value = "  +2.88%  "
$ws.Range("E27").Value = "  -3.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "156.84"
$ws.Range("E28").Value = "  -1.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "139.52"
$ws.Range("E29").Value = "  -4.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.287"
$ws.Range("E30").Value = "  -9.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.793"
$ws.Range("E31").Value = "  -6.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.490"
$ws.Range("E32").Value = "  -1.66%  "
$ws.Range("D33").Value = "1.794.82"
$ws.Range("E33").Value = "  -3.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08081"
$ws.Range("E34").Value = "  -2.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9776"
$ws.Range("E35").Value = "  -1.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02873"
$ws.Range("E36").Value = "  -6.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.552"
$ws.Range("E37").Value = "  -5.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2656"
$ws.Range("E38").Value = "  -5.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09112"
$ws.Range("E39").Value = "  -5.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.33"
$ws.Range("E40").Value = "  +0.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.55"
$ws.Range("E41").Value = "  +0.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.423"
$ws.Range("E42").Value = "  -5.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7471"
$ws.Range("E43").Value = "  -4.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.01"
$ws.Range("E44").Value = "  -3.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6890"
$ws.Range("E45").Value = "  -2.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.462"
$ws.Range("E46").Value = "  -3.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.066"
$ws.Range("E47").Value = "  -2.47%  "
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("E49").Value = "  -4.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.49"
$ws.Range("E50").Value = "  -2.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.199"
$ws.Range("E51").Value = "  -9.42%  "
